$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 4.922396194969107
$ws.Range("D2").Value = 3.908425251125884
$ws.Range("E2").Value = 16.5776391118889
$ws.Range("F2").Value = 18.79361713997934
$ws.Range("G2").Value = 19.38262887065268
$ws.Range("H2").Value = 11.86755819350586
$ws.Range("K2").Value = 15.35345112616489
$ws.Range("O2").Value = 16.80946291675029

$ws.Range("C3").Value = 4.749656062919869
$ws.Range("D3").Value = 3.832194808736343
$ws.Range("E3").Value = 15.63195195958803
$ws.Range("F3").Value = 18.86642199658159
$ws.Range("G3").Value = 19.48331001603888
$ws.Range("H3").Value = 11.94244192332396
$ws.Range("K3").Value = 14.52557911818607
$ws.Range("O3").Value = 16.93241433260239

$ws.Range("C4").Value = 4.641814755146623
$ws.Range("D4").Value = 3.784355124520419
$ws.Range("E4").Value = 15.02601058613874
$ws.Range("F4").Value = 18.9202297613242
$ws.Range("G4").Value = 19.5596728335035
$ws.Range("H4").Value = 11.99158451527076
$ws.Range("K4").Value = 13.99132957241778
$ws.Range("O4").Value = 17.01463168044188

$ws.Range("C5").Value = 4.597500425413741
$ws.Range("D5").Value = 3.764618242524349
$ws.Range("E5").Value = 14.77299341089933
$ws.Range("F5").Value = 18.9444259588477
$ws.Range("G5").Value = 19.59439493722446
$ws.Range("H5").Value = 12.01240414411071
$ws.Range("K5").Value = 13.76725935046791
$ws.Range("O5").Value = 17.0498145999936

$ws.Range("C6").Value = 4.590122267673529
$ws.Range("D6").Value = 3.761326889107455
$ws.Range("E6").Value = 14.73062106720035
$ws.Range("F6").Value = 18.94858014446695
$ws.Range("G6").Value = 19.60037654102408
$ws.Range("H6").Value = 12.01590910596154
$ws.Range("K6").Value = 13.7296739666152
$ws.Range("O6").Value = 17.05575772812934

$ws.Range("C7").Value = 4.641218497888531
$ws.Range("D7").Value = 3.784089900598135
$ws.Range("E7").Value = 15.02262257881794
$ws.Range("F7").Value = 18.92054692026751
$ws.Range("G7").Value = 19.56012659015786
$ws.Range("H7").Value = 11.99186208522967
$ws.Range("K7").Value = 13.98833319314668
$ws.Range("O7").Value = 17.01509938891241

$ws.Range("C8").Value = 4.86324982896815
$ws.Range("D8").Value = 3.882366215337472
$ws.Range("E8").Value = 16.25695879599966
$ws.Range("F8").Value = 18.81681849277812
$ws.Range("G8").Value = 19.41429381270279
$ws.Range("H8").Value = 11.89272031922256
$ws.Range("K8").Value = 15.07347047451003
$ws.Range("O8").Value = 16.85045343407776

$ws.Range("C9").Value = 5.281408412192639
$ws.Range("D9").Value = 4.066128010381672
$ws.Range("E9").Value = 18.57999496940994
$ws.Range("F9").Value = 18.68650825181048
$ws.Range("G9").Value = 19.24593930504486
$ws.Range("H9").Value = 11.72349316923153
$ws.Range("K9").Value = 16.99006086201405
$ws.Range("O9").Value = 16.58149472829069

$ws.Range("C10").Value = 5.574389549265597
$ws.Range("D10").Value = 4.194730873136172
$ws.Range("E10").Value = 20.23262884167457
$ws.Range("F10").Value = 18.63640362995793
$ws.Range("G10").Value = 19.19670165509524
$ws.Range("H10").Value = 11.61464056659879
$ws.Range("K10").Value = 18.2635693095583
$ws.Range("O10").Value = 16.41752142505423

$ws.Range("C11").Value = 5.703924806459257
$ws.Range("D11").Value = 4.251663005321444
$ws.Range("E11").Value = 20.94204666506026
$ws.Range("F11").Value = 18.62371482260991
$ws.Range("G11").Value = 19.19096042784507
$ws.Range("H11").Value = 11.56850965035772
$ws.Range("K11").Value = 18.81287239101498
$ws.Range("O11").Value = 16.35039442714788

$ws.Range("C12").Value = 5.752391626059548
$ws.Range("D12").Value = 4.272982084673106
$ws.Range("E12").Value = 21.20462853075963
$ws.Range("F12").Value = 18.62037627988761
$ws.Range("G12").Value = 19.19121446982395
$ws.Range("H12").Value = 11.55153060528103
$ws.Range("K12").Value = 19.01651290701051
$ws.Range("O12").Value = 16.326062562819

$ws.Range("C13").Value = 5.74198015372837
$ws.Range("D13").Value = 4.268401522876401
$ws.Range("E13").Value = 21.14834577623792
$ws.Range("F13").Value = 18.62102988016953
$ws.Range("G13").Value = 19.19105130359393
$ws.Range("H13").Value = 11.55516552102052
$ws.Range("K13").Value = 18.97285026606564
$ws.Range("O13").Value = 16.33125426434277

$ws.Range("C14").Value = 5.70792418223686
$ws.Range("D14").Value = 4.253421823864604
$ws.Range("E14").Value = 20.96377067594591
$ws.Range("F14").Value = 18.62341069865677
$ws.Range("G14").Value = 19.19093249695107
$ws.Range("H14").Value = 11.56710294318402
$ws.Range("K14").Value = 18.82971382459624
$ws.Range("O14").Value = 16.3483707500673

$ws.Range("C15").Value = 5.686986382887315
$ws.Range("D15").Value = 4.244214683255674
$ws.Range("E15").Value = 20.84992500742706
$ws.Range("F15").Value = 18.62506035967847
$ws.Range("G15").Value = 19.19117681541086
$ws.Range("H15").Value = 11.57447881906722
$ws.Range("K15").Value = 18.74146837618057
$ws.Range("O15").Value = 16.3589971605079

$ws.Range("C16").Value = 5.565844634954329
$ws.Range("D16").Value = 4.190977402995774
$ws.Range("E16").Value = 20.18541696523233
$ws.Range("F16").Value = 18.63743738786482
$ws.Range("G16").Value = 19.1974149780993
$ws.Range("H16").Value = 11.61772371054314
$ws.Range("K16").Value = 18.22706225127283
$ws.Range("O16").Value = 16.42205976315775

$ws.Range("C17").Value = 5.490534012945338
$ws.Range("D17").Value = 4.157905590900858
$ws.Range("E17").Value = 19.76693417399366
$ws.Range("F17").Value = 18.64762836514551
$ws.Range("G17").Value = 19.20553332810709
$ws.Range("H17").Value = 11.64512237558343
$ws.Range("K17").Value = 17.90376066712924
$ws.Range("O17").Value = 16.46266873133679

$ws.Range("C18").Value = 5.446867432288008
$ws.Range("D18").Value = 4.138736569861665
$ws.Range("E18").Value = 19.52224351761889
$ws.Range("F18").Value = 18.65444029404696
$ws.Range("G18").Value = 19.21176902369409
$ws.Range("H18").Value = 11.66119996421284
$ws.Range("K18").Value = 17.71498250492758
$ws.Range("O18").Value = 16.48672775094451

$ws.Range("C19").Value = 5.432024110754572
$ws.Range("D19").Value = 4.132221469622067
$ws.Range("E19").Value = 19.43870906662482
$ws.Range("F19").Value = 18.65690947801802
$ws.Range("G19").Value = 19.21414823393435
$ws.Range("H19").Value = 11.66669819781561
$ws.Range("K19").Value = 17.65058227416553
$ws.Range("O19").Value = 16.49499382816787

$ws.Range("C20").Value = 5.498587529978948
$ws.Range("D20").Value = 4.161441460027674
$ws.Range("E20").Value = 19.81189512841569
$ws.Range("F20").Value = 18.64644505703207
$ws.Range("G20").Value = 19.20450676066537
$ws.Range("H20").Value = 11.64217274735032
$ws.Range("K20").Value = 17.93846932091765
$ws.Range("O20").Value = 16.45827309954197

$ws.Range("C21").Value = 5.717943480681967
$ws.Range("D21").Value = 4.257828342667273
$ws.Range("E21").Value = 21.01814903395398
$ws.Range("F21").Value = 18.62267149786457
$ws.Range("G21").Value = 19.1909012612582
$ws.Range("H21").Value = 11.56358331953232
$ws.Range("K21").Value = 18.87187538082449
$ws.Range("O21").Value = 16.34331359094224

$ws.Range("C22").Value = 5.857875334320354
$ws.Range("D22").Value = 4.319418653726564
$ws.Range("E22").Value = 21.77120326455838
$ws.Range("F22").Value = 18.61568650165557
$ws.Range("G22").Value = 19.19617338890291
$ws.Range("H22").Value = 11.51507644810884
$ws.Range("K22").Value = 19.45643068331978
$ws.Range("O22").Value = 16.27452730476131

$ws.Range("C23").Value = 5.783519027968699
$ws.Range("D23").Value = 4.286679634064768
$ws.Range("E23").Value = 21.37250304613747
$ws.Range("F23").Value = 18.61862801754996
$ws.Range("G23").Value = 19.19205416939775
$ws.Range("H23").Value = 11.5407032409005
$ws.Range("K23").Value = 19.1467878561378
$ws.Range("O23").Value = 16.31065444912488

$ws.Range("C24").Value = 5.494947682240323
$ws.Range("D24").Value = 4.159843377138015
$ws.Range("E24").Value = 19.7915810527331
$ws.Range("F24").Value = 18.64697706304827
$ws.Range("G24").Value = 19.20496598895075
$ws.Range("H24").Value = 11.64350525977765
$ws.Range("K24").Value = 17.92278659161718
$ws.Range("O24").Value = 16.46025814660495

$ws.Range("C25").Value = 5.170542077602847
$ws.Range("D25").Value = 4.01747586668752
$ws.Range("E25").Value = 17.93359550080326
$ws.Range("F25").Value = 18.71381825512076
$ws.Range("G25").Value = 19.27859106995845
$ws.Range("H25").Value = 11.76656455851996
$ws.Range("K25").Value = 16.49487722653124
$ws.Range("O25").Value = 16.64840504627985
